# Auto-generated: apply scheduled market-data refresh to Chocobo_Profits sheets
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1900
$ws.Range("I62").Value = 1865.6666
$ws.Range("J62").Value = 2003
$ws.Range("K62").Value = 1865.6666
$ws.Range("L62").Value = 2003
$ws.Range("M62").Value = -1241.6666
$ws.Range("N62").Value = -3251
$ws.Range("H65").Value = 1900
$ws.Range("I65").Value = 1865.6666
$ws.Range("J65").Value = 2003
$ws.Range("K65").Value = 9328.333000000001
$ws.Range("L65").Value = 10015
$ws.Range("M65").Value = -6208.333000000001
$ws.Range("N65").Value = -16255
$ws.Range("H112").Value = 1370.5
$ws.Range("J112").Value = 1419.8246
$ws.Range("L112").Value = 4259.4738
$ws.Range("N112").Value = -6475.4738
$ws.Range("H135").Value = 1436.9412
$ws.Range("I135").Value = 673.4286
$ws.Range("K135").Value = 6060.8574
$ws.Range("M135").Value = -3525.8574
$ws.Range("H140").Value = 64537.066
$ws.Range("J140").Value = 64537.066
$ws.Range("L140").Value = 64537.066
$ws.Range("N140").Value = -74897.06599999999
$ws.Range("H141").Value = 9131.625
$ws.Range("I141").Value = 13620.667
$ws.Range("J141").Value = 3360
$ws.Range("K141").Value = 40862.001
$ws.Range("L141").Value = 10080
$ws.Range("M141").Value = -35682.001
$ws.Range("N141").Value = -20440

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4646.778
$ws.Range("I32").Value = 4753.173
$ws.Range("K32").Value = 4753.173
$ws.Range("M32").Value = -4466.173
$ws.Range("H61").Value = 1860.4
$ws.Range("I61").Value = 1825.5
$ws.Range("K61").Value = 1825.5
$ws.Range("M61").Value = -1613.5
$ws.Range("H110").Value = 2490.3635
$ws.Range("I110").Value = 2224.2
$ws.Range("J110").Value = 2712.1667
$ws.Range("K110").Value = 2224.2
$ws.Range("L110").Value = 2712.1667
$ws.Range("M110").Value = -179.1999999999998
$ws.Range("N110").Value = -6802.1667
$ws.Range("H122").Value = 2256.6667
$ws.Range("I122").Value = 1370.3334
$ws.Range("K122").Value = 4111.0002
$ws.Range("M122").Value = -1661.0002
$ws.Range("H132").Value = 3026.7942
$ws.Range("I132").Value = 2392.1428
$ws.Range("J132").Value = 3471.05
$ws.Range("K132").Value = 7176.428400000001
$ws.Range("L132").Value = 10413.15
$ws.Range("M132").Value = -4646.428400000001
$ws.Range("N132").Value = -15473.15
$ws.Range("H136").Value = 1860.4
$ws.Range("I136").Value = 1825.5
$ws.Range("K136").Value = 5476.5
$ws.Range("M136").Value = -2926.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1412.8182
$ws.Range("I107").Value = 1468.6428
$ws.Range("J107").Value = 1315.125
$ws.Range("K107").Value = 1468.6428
$ws.Range("L107").Value = 1315.125
$ws.Range("M107").Value = 451.3571999999999
$ws.Range("N107").Value = -5155.125
$ws.Range("H134").Value = 2722.775
$ws.Range("I134").Value = 1137.579
$ws.Range("J134").Value = 4157
$ws.Range("K134").Value = 3412.737
$ws.Range("L134").Value = 12471
$ws.Range("M134").Value = -877.7370000000001
$ws.Range("N134").Value = -17541

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 29930
$ws.Range("J9").Value = 29930
$ws.Range("L9").Value = 29930
$ws.Range("N9").Value = -30266
$ws.Range("H99").Value = 12504227
$ws.Range("I99").Value = 25002054
$ws.Range("J99").Value = 6400.375
$ws.Range("K99").Value = 25002054
$ws.Range("L99").Value = 6400.375
$ws.Range("M99").Value = -25000556
$ws.Range("N99").Value = -9396.375
$ws.Range("H126").Value = 12504227
$ws.Range("I126").Value = 25002054
$ws.Range("J126").Value = 6400.375
$ws.Range("K126").Value = 75006162
$ws.Range("L126").Value = 19201.125
$ws.Range("M126").Value = -75003692
$ws.Range("N126").Value = -24141.125
$ws.Range("H132").Value = 5150.222
$ws.Range("I132").Value = 4157.6665
$ws.Range("K132").Value = 12472.9995
$ws.Range("M132").Value = -9942.999500000002
$ws.Range("H137").Value = 48780
$ws.Range("J137").Value = 48780
$ws.Range("L137").Value = 48780
$ws.Range("N137").Value = -58980

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3168.3818
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 3348.1836
$ws.Range("K68").Value = 5100
$ws.Range("L68").Value = 10044.5508
$ws.Range("M68").Value = -4289
$ws.Range("N68").Value = -11666.5508
$ws.Range("H71").Value = 3168.3818
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 3348.1836
$ws.Range("K71").Value = 15300
$ws.Range("L71").Value = 30133.6524
$ws.Range("M71").Value = -11244
$ws.Range("N71").Value = -38245.6524
$ws.Range("H131").Value = 772.9400000000001
$ws.Range("I131").Value = 451.33334
$ws.Range("J131").Value = 804.74725
$ws.Range("K131").Value = 1354.00002
$ws.Range("L131").Value = 2414.24175
$ws.Range("M131").Value = 3685.99998
$ws.Range("N131").Value = -12494.24175

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6390.659
$ws.Range("I70").Value = 5830
$ws.Range("J70").Value = 8913.625
$ws.Range("K70").Value = 5830
$ws.Range("L70").Value = 8913.625
$ws.Range("M70").Value = -5560
$ws.Range("N70").Value = -9453.625
$ws.Range("H73").Value = 6390.659
$ws.Range("I73").Value = 5830
$ws.Range("J73").Value = 8913.625
$ws.Range("K73").Value = 5830
$ws.Range("L73").Value = 8913.625
$ws.Range("M73").Value = -4894
$ws.Range("N73").Value = -10785.625
$ws.Range("H102").Value = 2256.258
$ws.Range("I102").Value = 1823.8518
$ws.Range("J102").Value = 5175
$ws.Range("K102").Value = 1823.8518
$ws.Range("L102").Value = 5175
$ws.Range("M102").Value = -201.8517999999999
$ws.Range("N102").Value = -8419
$ws.Range("H126").Value = 3665.5212
$ws.Range("I126").Value = 2851.16
$ws.Range("J126").Value = 5604.476
$ws.Range("K126").Value = 8553.48
$ws.Range("L126").Value = 16813.428
$ws.Range("M126").Value = -6083.48
$ws.Range("N126").Value = -21753.428
$ws.Range("H132").Value = 4513.8486
$ws.Range("I132").Value = 3513.9443
$ws.Range("K132").Value = 10541.8329
$ws.Range("M132").Value = -8011.832900000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2966.3076
$ws.Range("I7").Value = 1506.5555
$ws.Range("J7").Value = 6250.75
$ws.Range("K7").Value = 1506.5555
$ws.Range("L7").Value = 6250.75
$ws.Range("M7").Value = -1394.5555
$ws.Range("N7").Value = -6474.75
$ws.Range("H40").Value = 4811.567
$ws.Range("I40").Value = 4047.3157
$ws.Range("J40").Value = 6131.636
$ws.Range("K40").Value = 4047.3157
$ws.Range("L40").Value = 6131.636
$ws.Range("M40").Value = -3911.3157
$ws.Range("N40").Value = -6403.636
$ws.Range("H122").Value = 3794.913
$ws.Range("I122").Value = 2892.8235
$ws.Range("J122").Value = 6350.8335
$ws.Range("K122").Value = 8678.470499999999
$ws.Range("L122").Value = 19052.5005
$ws.Range("M122").Value = -6228.470499999999
$ws.Range("N122").Value = -23952.5005
$ws.Range("H126").Value = 2966.3076
$ws.Range("I126").Value = 1506.5555
$ws.Range("J126").Value = 6250.75
$ws.Range("K126").Value = 4519.666499999999
$ws.Range("L126").Value = 18752.25
$ws.Range("M126").Value = -2049.666499999999
$ws.Range("N126").Value = -23692.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3028.7036
$ws.Range("I122").Value = 886.5625
$ws.Range("K122").Value = 2659.6875
$ws.Range("M122").Value = -209.6875
